$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 274; this shifts the existing rows 274-374
# down to 275-375 and carries formatting (e.g. the date number format on
# column D) down from the row above, matching native Excel behaviour.
$ws.Rows("274:274").Insert()

# Populate the newly inserted row 274 with the new data record.
$ws.Range("A274").Value = 8
$ws.Range("B274").Value = "Terminal La Palmera de La Serena"
$ws.Range("C274").Value = "Coquimbo"
$ws.Range("D274").Value = 45007
$ws.Range("E274").Value = 4
$ws.Range("F274").Value = "Fruta"
$ws.Range("G274").Value = 100103
$ws.Range("H274").Value = "Frutos de hueso (carozo)"
$ws.Range("I274").Value = 100103002
$ws.Range("J274").Value = "Ciruela"
$ws.Range("K274").Value = "Angeleno"
$ws.Range("L274").Value = "Primera"
$ws.Range("M274").Value = 22
$ws.Range("N274").Value = 205000
$ws.Range("O274").Value = 210000
$ws.Range("P274").Value = 207500
$ws.Range("Q274").Value = "`$/bins (450 kilos)"
$ws.Range("R274").Value = "Región Metropolitana"
$ws.Range("S274").Value = 461
$ws.Range("T274").Value = 450
